$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# Row 7 ("Experimental"): set the Value column (B7) to the text "false".
# Typing the bare word false/true into a cell gets auto-recognized as a
# Boolean rather than text, so it is built up as a formula that evaluates
# to the text string, then converted to a plain value in place (this
# keeps it as a real text cell without disturbing the cell's style).
$wsMeta.Range("B7").Formula = "=TRIM("" false"")"
$wsMeta.Range("B7").Copy()
$wsMeta.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Row 8 ("Date"): update the Value column (B8) to the new timestamp
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
